# Updated sample collection method id and context instructions
#
# The "Instructions" tab's row describing the Sampling Method Context now
# also tells users to enter the context "MassWateR" when using MassWateR's
# standard sampling methods.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Meta")
$wsInstructions = $wb.Worksheets.Item("Instructions")

# Update the instructional text for the "Sampling Method Context" row (B7)
# on the Instructions tab.
$wsInstructions.Range("B7").Value = 'Enter the Context for the Sampling Method IDs that are used for sampling this parameter.  Not applicable for field measurements/observations.  If you are using the standard methods defined by MassWateR, enter the context "MassWateR".'

# Reflect the cell that was left selected on the Instructions tab after the
# edit, then restore the Meta tab as the active sheet (it was active before
# this edit and should remain so).
$wsInstructions.Activate()
$wsInstructions.Range("B8").Select()
$wsMeta.Activate()
